$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (x values) per diff
$ws.Range("A3:A7").Value = 4.5
$ws.Range("A8:A12").Value = 3.6
$ws.Range("A13:A17").Value = 2.7
$ws.Range("A18:A22").Value = 1.8

# Update the view: scroll so row 13 is the top-left visible row,
# and set the active selection to D9
$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D9").Select()
